$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FirstSet")

# Update the user emails referenced in column A for several rows.
$ws.Range("A2").Value = "qa@cvhcare.com"
$ws.Range("A8").Value = "qa@cvhcare.com"
$ws.Range("A9").Value = "qa@cvhcare.com"
$ws.Range("A20").Value = "kavya.mothukuri@cvhcare.com"
$ws.Range("A22").Value = "russell.sadang@cvhcare.com"

# Row 2 previously stored the literal password "Password0!" in column B; normalize it
# back to the generic "password" placeholder used by every other row.
$ws.Range("B2").Value = "password"

# The hyperlink that used to live on A8 (mailto:scheduler.user1@cvhcare.com) needs to
# be removed from its original position and re-created at the end of the hyperlinks
# collection, now pointing at marc.miller@cvhcare.com (mirrors what Excel does when a
# row above it is removed and the sheet's hyperlinks are left pointing at stale cells).
$existing = $null
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$8') {
        $existing = $hl
    }
}
if ($existing -ne $null) {
    $existing.Delete()
}
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:marc.miller@cvhcare.com")

# Re-creating the hyperlink re-stamps A8 with a fresh "visited hyperlink" style; put it
# back on the same Hyperlink style used by every other linked cell in the column.
$ws.Range("A8").Style = $ws.Range("A9").Style

# Update the active selection to match the new state of the sheet.
$ws.Range("A22").Select()
